$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value2 = 6428369
$ws.Range('F2').Value2 = 'Nykobing'
$ws.Range('G2').Value2 = 'FC Fredericia'
$ws.Range('H2').Value2 = 1
$ws.Range('I2').Value2 = 5
$ws.Range('J2').Value2 = 'A'
$ws.Range('K2').Value2 = 4.2
$ws.Range('L2').Value2 = 4.2
$ws.Range('M2').Value2 = 1.65
$ws.Range('N2').Value2 = 6
$ws.Range('O2').Value2 = 5.75
$ws.Range('P2').Value2 = 1.444
$ws.Range('Q2').Value2 = 1.5
$ws.Range('R2').Value2 = 1.825
$ws.Range('S2').Value2 = 2.025
$ws.Range('T2').Value2 = 4
$ws.Range('U2').Value2 = 1.875
$ws.Range('V2').Value2 = 1.975
$ws.Range('W2').Value2 = -1
$ws.Range('Y2').Value2 = 0.444
$ws.Range('Z2').Value2 = -1
$ws.Range('AA2').Value2 = 1.025
$ws.Range('AB2').Value2 = 0.875
$ws.Range('B4').Value2 = 6428368
$ws.Range('F4').Value2 = 'HB Kge'
$ws.Range('G4').Value2 = 'Fremad Amager'
$ws.Range('H4').Value2 = 4
$ws.Range('I4').Value2 = 2
$ws.Range('J4').Value2 = 'H'
$ws.Range('K4').Value2 = 1.571
$ws.Range('L4').Value2 = 4.333
$ws.Range('M4').Value2 = 4.333
$ws.Range('N4').Value2 = 1.45
$ws.Range('O4').Value2 = 5
$ws.Range('P4').Value2 = 6
$ws.Range('Q4').Value2 = -1.25
$ws.Range('R4').Value2 = 1.925
$ws.Range('S4').Value2 = 1.925
$ws.Range('T4').Value2 = 3.5
$ws.Range('U4').Value2 = 1.9
$ws.Range('V4').Value2 = 1.95
$ws.Range('W4').Value2 = 0.45
$ws.Range('Y4').Value2 = -1
$ws.Range('Z4').Value2 = 0.925
$ws.Range('AA4').Value2 = -1
$ws.Range('AB4').Value2 = 0.8999999999999999
$ws.Range('B5').Value2 = 6428339
$ws.Range('F5').Value2 = 'Vejle'
$ws.Range('G5').Value2 = 'Vendsyssel FF'
$ws.Range('H5').Value2 = 4
$ws.Range('J5').Value2 = 'H'
$ws.Range('K5').Value2 = 1.55
$ws.Range('M5').Value2 = 5
$ws.Range('N5').Value2 = 2.2
$ws.Range('O5').Value2 = 3.8
$ws.Range('P5').Value2 = 3
$ws.Range('Q5').Value2 = -0.25
$ws.Range('R5').Value2 = 1.875
$ws.Range('S5').Value2 = 1.975
$ws.Range('T5').Value2 = 3
$ws.Range('U5').Value2 = 2.05
$ws.Range('V5').Value2 = 1.8
$ws.Range('W5').Value2 = 1.2
$ws.Range('Y5').Value2 = -1
$ws.Range('Z5').Value2 = 0.875
$ws.Range('AA5').Value2 = -1
$ws.Range('AB5').Value2 = 1.05
$ws.Range('B7').Value2 = 6428337
$ws.Range('F7').Value2 = 'Hvidovre IF'
$ws.Range('G7').Value2 = 'Naestved'
$ws.Range('H7').Value2 = 2
$ws.Range('J7').Value2 = 'A'
$ws.Range('K7').Value2 = 1.7
$ws.Range('M7').Value2 = 4.2
$ws.Range('N7').Value2 = 1.571
$ws.Range('O7').Value2 = 4.75
$ws.Range('P7').Value2 = 4.75
$ws.Range('Q7').Value2 = -1
$ws.Range('R7').Value2 = 1.975
$ws.Range('S7').Value2 = 1.875
$ws.Range('T7').Value2 = 3.5
$ws.Range('U7').Value2 = 2
$ws.Range('V7').Value2 = 1.85
$ws.Range('W7').Value2 = -1
$ws.Range('Y7').Value2 = 3.75
$ws.Range('Z7').Value2 = -1
$ws.Range('AA7').Value2 = 0.875
$ws.Range('AB7').Value2 = 1
$ws.Range('B17').Value2 = 6798549
$ws.Range('F17').Value2 = 'Kolding IF'
$ws.Range('G17').Value2 = 'Hillerd'
$ws.Range('H17').Value2 = 2
$ws.Range('I17').Value2 = 1
$ws.Range('K17').Value2 = 1.909
$ws.Range('L17').Value2 = 3.5
$ws.Range('M17').Value2 = 3.6
$ws.Range('N17').Value2 = 1.5
$ws.Range('P17').Value2 = 6
$ws.Range('R17').Value2 = 1.85
$ws.Range('S17').Value2 = 2
$ws.Range('U17').Value2 = 1.925
$ws.Range('V17').Value2 = 1.925
$ws.Range('W17').Value2 = 0.5
$ws.Range('AB17').Value2 = 0.4625
$ws.Range('AC17').Value2 = -0.5
$ws.Range('B18').Value2 = 6798865
$ws.Range('F18').Value2 = 'AaB'
$ws.Range('G18').Value2 = 'AC Horsens'
$ws.Range('H18').Value2 = 1
$ws.Range('I18').Value2 = 0
$ws.Range('K18').Value2 = 1.7
$ws.Range('L18').Value2 = 4
$ws.Range('M18').Value2 = 4.333
$ws.Range('N18').Value2 = 1.6
$ws.Range('P18').Value2 = 5
$ws.Range('R18').Value2 = 2.025
$ws.Range('S18').Value2 = 1.825
$ws.Range('U18').Value2 = 1.825
$ws.Range('V18').Value2 = 2.025
$ws.Range('W18').Value2 = 0.6000000000000001
$ws.Range('AB18').Value2 = -1
$ws.Range('AC18').Value2 = 1.025
$ws.Range('B105').Value2 = 6799300
$ws.Range('F105').Value2 = 'Sonderjyske'
$ws.Range('G105').Value2 = 'Naestved'
$ws.Range('H105').Value2 = 4
$ws.Range('I105').Value2 = 1
$ws.Range('J105').Value2 = 'H'
$ws.Range('K105').Value2 = 1.25
$ws.Range('L105').Value2 = 6
$ws.Range('M105').Value2 = 8
$ws.Range('N105').Value2 = 1.3
$ws.Range('O105').Value2 = 6
$ws.Range('P105').Value2 = 8.5
$ws.Range('Q105').Value2 = -1.5
$ws.Range('R105').Value2 = 1.875
$ws.Range('S105').Value2 = 1.975
$ws.Range('T105').Value2 = 3.25
$ws.Range('U105').Value2 = 1.925
$ws.Range('V105').Value2 = 1.925
$ws.Range('W105').Value2 = 0.3
$ws.Range('Y105').Value2 = -1
$ws.Range('Z105').Value2 = 0.875
$ws.Range('AA105').Value2 = -1
$ws.Range('AB105').Value2 = 0.925
$ws.Range('B106').Value2 = 6799298
$ws.Range('F106').Value2 = 'FC Helsingor'
$ws.Range('G106').Value2 = 'AC Horsens'
$ws.Range('H106').Value2 = 1
$ws.Range('I106').Value2 = 2
$ws.Range('J106').Value2 = 'A'
$ws.Range('K106').Value2 = 3.1
$ws.Range('L106').Value2 = 3.6
$ws.Range('M106').Value2 = 2.1
$ws.Range('N106').Value2 = 3.2
$ws.Range('O106').Value2 = 3.4
$ws.Range('P106').Value2 = 2.25
$ws.Range('Q106').Value2 = 0.25
$ws.Range('R106').Value2 = 1.95
$ws.Range('S106').Value2 = 1.9
$ws.Range('T106').Value2 = 2.5
$ws.Range('U106').Value2 = 1.85
$ws.Range('V106').Value2 = 2
$ws.Range('W106').Value2 = -1
$ws.Range('Y106').Value2 = 1.25
$ws.Range('Z106').Value2 = -1
$ws.Range('AA106').Value2 = 0.8999999999999999
$ws.Range('AB106').Value2 = 0.8500000000000001
$ws.Range('B125').Value2 = 6800815
$ws.Range('F125').Value2 = 'B93 Copenhagen'
$ws.Range('G125').Value2 = 'HB Kge'
$ws.Range('H125').Value2 = 3
$ws.Range('I125').Value2 = 1
$ws.Range('J125').Value2 = 'H'
$ws.Range('K125').Value2 = 2.375
$ws.Range('L125').Value2 = 3.4
$ws.Range('M125').Value2 = 2.75
$ws.Range('N125').Value2 = 2.3
$ws.Range('O125').Value2 = 3.4
$ws.Range('P125').Value2 = 3.1
$ws.Range('Q125').Value2 = -0.25
$ws.Range('T125').Value2 = 3
$ws.Range('U125').Value2 = 2.025
$ws.Range('V125').Value2 = 1.825
$ws.Range('W125').Value2 = 1.3
$ws.Range('Y125').Value2 = -1
$ws.Range('Z125').Value2 = 1.025
$ws.Range('AA125').Value2 = -1
$ws.Range('AB125').Value2 = 1.025
$ws.Range('AC125').Value2 = -1
$ws.Range('B126').Value2 = 6799308
$ws.Range('F126').Value2 = 'Sonderjyske'
$ws.Range('G126').Value2 = 'Vendsyssel FF'
$ws.Range('H126').Value2 = 0
$ws.Range('I126').Value2 = 3
$ws.Range('J126').Value2 = 'A'
$ws.Range('K126').Value2 = 1.666
$ws.Range('L126').Value2 = 4
$ws.Range('M126').Value2 = 4.5
$ws.Range('N126').Value2 = 1.8
$ws.Range('O126').Value2 = 3.75
$ws.Range('P126').Value2 = 4.5
$ws.Range('Q126').Value2 = -0.75
$ws.Range('T126').Value2 = 2.75
$ws.Range('U126').Value2 = 1.85
$ws.Range('V126').Value2 = 2
$ws.Range('W126').Value2 = -1
$ws.Range('Y126').Value2 = 3.5
$ws.Range('Z126').Value2 = -1
$ws.Range('AA126').Value2 = 0.825
$ws.Range('AB126').Value2 = 0.425
$ws.Range('AC126').Value2 = -0.5
$ws.Range('N152').Value2 = 2.4
$ws.Range('N153').Value2 = 2.3
$ws.Range('O153').Value2 = 3.5
$ws.Range('N154').Value2 = 5
$ws.Range('O154').Value2 = 3.75
$ws.Range('P154').Value2 = 1.727
$ws.Range('R154').Value2 = 1.875
$ws.Range('S154').Value2 = 1.975
$ws.Range('P155').Value2 = 3.4
$ws.Range('N156').Value2 = 4.333
$ws.Range('N157').Value2 = 4
$ws.Range('O157').Value2 = 3.8
$ws.Range('P157').Value2 = 1.85
